$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2249.5
$ws.Range("I19").Value = 2199.5
$ws.Range("K19").Value = 2199.5
$ws.Range("M19").Value = -2024.5
$ws.Range("H32").Value = 4144.3335
$ws.Range("I32").Value = 2049
$ws.Range("K32").Value = 2049
$ws.Range("M32").Value = -1723
$ws.Range("H96").Value = 2080
$ws.Range("I96").Value = 3400
$ws.Range("J96").Value = 760
$ws.Range("K96").Value = 10200
$ws.Range("L96").Value = 2280
$ws.Range("M96").Value = -8827
$ws.Range("N96").Value = -5026
$ws.Range("H127").Value = 2305.4443
$ws.Range("I127").Value = 2011.8572
$ws.Range("J127").Value = 3333
$ws.Range("K127").Value = 6035.571599999999
$ws.Range("L127").Value = 9999
$ws.Range("M127").Value = -1075.571599999999
$ws.Range("N127").Value = -19919
$ws.Range("H129").Value = 2269.625
$ws.Range("I129").Value = 1093.2858
$ws.Range("J129").Value = 3184.5557
$ws.Range("K129").Value = 3279.8574
$ws.Range("L129").Value = 9553.667099999999
$ws.Range("M129").Value = 1720.1426
$ws.Range("N129").Value = -19553.6671
$ws.Range("H137").Value = 1122.25
$ws.Range("J137").Value = 1096.6666
$ws.Range("L137").Value = 3289.9998
$ws.Range("N137").Value = -8389.9998
$ws.Range("H138").Value = 16461.47
$ws.Range("J138").Value = 19846.285
$ws.Range("L138").Value = 59538.855
$ws.Range("N138").Value = -69818.855
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 11835168
$ws.Range("I13").Value = 11835168
$ws.Range("K13").Value = 11835168
$ws.Range("M13").Value = -11835024
$ws.Range("H80").Value = 57498.5
$ws.Range("J80").Value = 69998.336
$ws.Range("L80").Value = 69998.336
$ws.Range("N80").Value = -71994.336
$ws.Range("H83").Value = 57498.5
$ws.Range("J83").Value = 69998.336
$ws.Range("L83").Value = 209995.008
$ws.Range("N83").Value = -219979.008
$ws.Range("H122").Value = 1398.5454
$ws.Range("I122").Value = 1376
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 4128
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -1678
$ws.Range("N122").Value = -9400
$ws.Range("H124").Value = 65857
$ws.Range("J124").Value = 65857
$ws.Range("L124").Value = 65857
$ws.Range("N124").Value = -75677
$ws.Range("H132").Value = 3500.8667
$ws.Range("I132").Value = 3500.8667
$ws.Range("K132").Value = 10502.6001
$ws.Range("M132").Value = -7972.6001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 24882
$ws.Range("I96").Value = 24882
$ws.Range("K96").Value = 24882
$ws.Range("M96").Value = -22136
$ws.Range("H99").Value = 1183.65
$ws.Range("J99").Value = 1088.579
$ws.Range("L99").Value = 1088.579
$ws.Range("N99").Value = -4084.579
$ws.Range("H105").Value = 4936.6665
$ws.Range("I105").Value = 4936.6665
$ws.Range("K105").Value = 4936.6665
$ws.Range("M105").Value = -3189.6665
$ws.Range("H134").Value = 2063.6667
$ws.Range("I134").Value = 2074.182
$ws.Range("J134").Value = 1948
$ws.Range("K134").Value = 6222.545999999999
$ws.Range("L134").Value = 5844
$ws.Range("M134").Value = -3687.545999999999
$ws.Range("N134").Value = -10914
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1566.3334
$ws.Range("I16").Value = 1566.3334
$ws.Range("K16").Value = 1566.3334
$ws.Range("M16").Value = -1279.3334
$ws.Range("H113").Value = 1566.3334
$ws.Range("I113").Value = 1566.3334
$ws.Range("K113").Value = 1566.3334
$ws.Range("M113").Value = 603.6666
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
$ws.Range("H141").Value = 328090.8
$ws.Range("J141").Value = 496484.66
$ws.Range("L141").Value = 496484.66
$ws.Range("N141").Value = -506844.66
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 613.125
$ws.Range("I14").Value = 613.125
$ws.Range("K14").Value = 1839.375
$ws.Range("M14").Value = -1666.375
$ws.Range("H97").Value = 1723
$ws.Range("I97").Value = 1723
$ws.Range("K97").Value = 5169
$ws.Range("M97").Value = -4673
$ws.Range("H99").Value = 791
$ws.Range("I99").Value = 1045.3334
$ws.Range("K99").Value = 3136.0002
$ws.Range("M99").Value = -890.0001999999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9249.5
$ws.Range("I70").Value = 7999
$ws.Range("J70").Value = 10500
$ws.Range("K70").Value = 7999
$ws.Range("L70").Value = 10500
$ws.Range("M70").Value = -7729
$ws.Range("N70").Value = -11040
$ws.Range("H73").Value = 9249.5
$ws.Range("I73").Value = 7999
$ws.Range("J73").Value = 10500
$ws.Range("K73").Value = 7999
$ws.Range("L73").Value = 10500
$ws.Range("M73").Value = -7063
$ws.Range("N73").Value = -12372
$ws.Range("H113").Value = 2805.625
$ws.Range("I113").Value = 2744.3333
$ws.Range("J113").Value = 2989.5
$ws.Range("K113").Value = 2744.3333
$ws.Range("L113").Value = 2989.5
$ws.Range("M113").Value = -574.3332999999998
$ws.Range("N113").Value = -7329.5
$ws.Range("H122").Value = 4197.6
$ws.Range("I122").Value = 3749.25
$ws.Range("J122").Value = 4496.5
$ws.Range("K122").Value = 11247.75
$ws.Range("L122").Value = 13489.5
$ws.Range("M122").Value = -8797.75
$ws.Range("N122").Value = -18389.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 811.25
$ws.Range("I16").Value = 811.25
$ws.Range("K16").Value = 811.25
$ws.Range("M16").Value = -641.25
$ws.Range("H122").Value = 3488.2666
$ws.Range("J122").Value = 3476.348
$ws.Range("L122").Value = 10429.044
$ws.Range("N122").Value = -15329.044
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1639.4546
$ws.Range("J14").Value = 4252.5
$ws.Range("L14").Value = 4252.5
$ws.Range("N14").Value = -4588.5
$ws.Range("H81").Value = 1464.4117
$ws.Range("I81").Value = 1593
$ws.Range("J81").Value = 500
$ws.Range("K81").Value = 3186
$ws.Range("L81").Value = 1000
$ws.Range("M81").Value = -2125
$ws.Range("N81").Value = -3122
$ws.Range("H84").Value = 1464.4117
$ws.Range("I84").Value = 1593
$ws.Range("J84").Value = 500
$ws.Range("K84").Value = 15930
$ws.Range("L84").Value = 5000
$ws.Range("M84").Value = -10626
$ws.Range("N84").Value = -15608
